$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 46
$ws.Range("F12").Value = 14.5

$ws.Range("E21").Value = 4.399999999999999

$ws.Range("A27").Value = 'Örebro University'
$ws.Range("B27").Value = 17
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 18.4

$ws.Range("A28").Value = 'Örebro University Hospital'
$ws.Range("B28").Value = 1
$ws.Range("F28").Value = 94.89999999999999

$ws.Range("A29").Value = 'Oslo University Hospital'
$ws.Range("B29").Value = 102
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0.1
$ws.Range("F29").Value = 5.3

$ws.Range("A30").Value = 'Oulu University Hospital'
$ws.Range("B30").Value = 10
$ws.Range("F30").Value = 27.8

$ws.Range("A31").Value = 'Sahlgrenska University Hospital'
$ws.Range("B31").Value = 40
$ws.Range("F31").Value = 8.799999999999999

$ws.Range("A32").Value = 'Skane University Hospital'
$ws.Range("B32").Value = 23
$ws.Range("F32").Value = 14.3

$ws.Range("A33").Value = 'St. Olav’s University Hospital'
$ws.Range("B33").Value = 24
$ws.Range("F33").Value = 13.8

$ws.Range("A34").Value = 'Steno Diabetes Center Copenhagen'
$ws.Range("B34").Value = 13
$ws.Range("F34").Value = 22.8

$ws.Range("A35").Value = 'Stockholm South General Hospital'
$ws.Range("B35").Value = 3
$ws.Range("F35").Value = 56.10000000000001

$ws.Range("A36").Value = 'Tampere University Hospital'
$ws.Range("B36").Value = 22
$ws.Range("F36").Value = 14.9

$ws.Range("A37").Value = 'The National University Hospital of Iceland'
$ws.Range("B37").Value = 5
$ws.Range("F37").Value = 43.4

$ws.Range("A38").Value = 'Turku University Hospital'
$ws.Range("B38").Value = 48
$ws.Range("F38").Value = 7.399999999999999

$ws.Range("A39").Value = 'UiT The Arctic University of Norway'
$ws.Range("B39").Value = 14
$ws.Range("F39").Value = 21.5

$ws.Range("A40").Value = 'Umeå University'
$ws.Range("B40").Value = 41
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 8.6

$ws.Range("A41").Value = 'University Hospital of North Norway'
$ws.Range("B41").Value = 17
$ws.Range("F41").Value = 18.4

$ws.Range("A42").Value = 'University Hospital of Umeå'
$ws.Range("B42").Value = 2
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = 50
$ws.Range("E42").Value = 2.6
$ws.Range("F42").Value = 97.39999999999999

$ws.Range("A43").Value = 'University of Bergen'
$ws.Range("B43").Value = 31
$ws.Range("F43").Value = 11

$ws.Range("A44").Value = 'University of Copenhagen'
$ws.Range("B44").Value = 97
$ws.Range("F44").Value = 3.8

$ws.Range("A45").Value = 'University of Eastern Finland'
$ws.Range("B45").Value = 12
$ws.Range("F45").Value = 24.2

$ws.Range("A46").Value = 'University of Helsinki'
$ws.Range("B46").Value = 21
$ws.Range("F46").Value = 15.5

$ws.Range("A47").Value = 'University of Iceland'
$ws.Range("B47").Value = 5
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 43.4

$ws.Range("A48").Value = 'University of Oslo'
$ws.Range("B48").Value = 23
$ws.Range("F48").Value = 14.3

$ws.Range("A49").Value = 'University of Oulu'
$ws.Range("B49").Value = 25
$ws.Range("D49").Value = 4
$ws.Range("E49").Value = 0.2
$ws.Range("F49").Value = 19.5

$ws.Range("A50").Value = 'University of Southern Denmark'
$ws.Range("B50").Value = 42
$ws.Range("F50").Value = 8.4

$ws.Range("A51").Value = 'University of Tampere'
$ws.Range("C51").Value = 1
$ws.Range("D51").Value = 11.1
$ws.Range("E51").Value = 0.6
$ws.Range("F51").Value = 43.5

$ws.Range("A52").Value = 'University of Turku'
$ws.Range("B52").Value = 20
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 0
$ws.Range("F52").Value = 16.1

$ws.Range("A53").Value = 'Uppsala Academic Hospital'
$ws.Range("B53").Value = 9
$ws.Range("F53").Value = 29.9

$ws.Range("A54").Value = 'Uppsala University'
$ws.Range("B54").Value = 51
$ws.Range("C54").Value = 1
$ws.Range("D54").Value = 2
$ws.Range("E54").Value = 0.1
$ws.Range("F54").Value = 10.3

$ws.Range("A55").Value = 'Zealand University Hospital'
$ws.Range("B55").Value = 28
$ws.Range("F55").Value = 12.1
